# Daily attendance processing - 2026-01-08 04:31:11
# Reorders the "Recorded By" (column G) values for certain rows: the
# last name/email in the comma-separated list is moved to the front.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Exact mapping of old -> new "Recorded By" text observed in the edit.
$map = @{
    "System, system, backup@backdoor.com" = "backup@backdoor.com, System, system"
    "dnasr281@gmail.com, System"          = "System, dnasr281@gmail.com"
    "System, backup@backdoor.com"         = "backup@backdoor.com, System"
    "dnasr281@gmail.com, admin@admin.com" = "admin@admin.com, dnasr281@gmail.com"
}

# Find the last used row in column G (data starts at row 2; header at row 1).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row
if ($lastRow -lt 2) {
    $lastRow = $ws.UsedRange.Rows.Count
}

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $current = $cell.Value2
    if ($null -ne $current -and $map.ContainsKey($current)) {
        $cell.Value = $map[$current]
    }
}
